# Add two new columns "I0" (I) and "IF" (J) after the existing "IP" column (H).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1) -------------------------------------------------
# Copy the formatting of the existing header cell (H1, style "1": bold,
# bordered, centered) onto the two new header cells, then set their text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (2-27) ------------------------------------------------------
# Column I = starting index (I0), column J = final index (IF) for each row.
$i0 = @{
    2 = 1; 3 = 1; 4 = 1; 5 = 1; 6 = 1; 7 = 1; 8 = 1; 9 = 1; 10 = 1;
    11 = 1; 12 = 1; 13 = 1; 14 = 1; 15 = 1; 16 = 1; 17 = 1; 18 = 1;
    19 = 1; 20 = 1; 21 = 1; 22 = 7; 23 = 6; 24 = 1; 25 = 1; 26 = 1; 27 = 3
}
$if = @{
    2 = 4; 3 = 6; 4 = 6; 5 = 3; 6 = 8; 7 = 6; 8 = 5; 9 = 4; 10 = 9;
    11 = 6; 12 = 6; 13 = 6; 14 = 8; 15 = 6; 16 = 5; 17 = 6; 18 = 2;
    19 = 6; 20 = 5; 21 = 4; 22 = 9; 23 = 8; 24 = 6; 25 = 4; 26 = 3; 27 = 4
}

for ($row = 2; $row -le 27; $row++) {
    $ws.Cells.Item($row, 9).Value = $i0[$row]
    $ws.Cells.Item($row, 10).Value = $if[$row]
}
